# Auto-generated edit script: apply updated cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.817.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.513.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.30%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.57%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.513.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.573'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '

# Row 11
$ws.Range("E11").Value = '  +4.78%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.438'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.91%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.124.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.45%  '

# Row 14
$ws.Range("E14").Value = '  +0.03%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.15%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.785.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.18%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000178'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.516.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.51%  '

# Row 19
$ws.Range("E19").Value = '  +3.32%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.70%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.36%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.82%  '

# Row 24
$ws.Range("E24").Value = '  +0.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.99%  '

# Row 26
$ws.Range("E26").Value = '  +3.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.07%  '

# Row 28
$ws.Range("E28").Value = '  +0.96%  '

# Row 29
$ws.Range("E29").Value = '  +0.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.49%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.31%  '

# Row 32
$ws.Range("E32").Value = '  +3.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.08%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.50%  '

# Row 35
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("E36").Value = '  +5.54%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.899'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.74%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0747'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.73%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.76%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.840.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.27%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.90%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.50%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0314'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.81%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '352.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.28%  '

# Row 50
$ws.Range("E50").Value = '  +2.30%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.51%  '
